# Insert a new daily price record at the top of the "Ajo" weekly block
# (row 309), pushing the existing rows 309-335 down to 310-336.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(309).Insert()

$ws.Cells.Item(309, 1).Value = 11
$ws.Cells.Item(309, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(309, 3).Value = "Bíobío"
$ws.Cells.Item(309, 4).Value = 45212
$ws.Cells.Item(309, 5).Value = 8
$ws.Cells.Item(309, 6).Value = 100112003
$ws.Cells.Item(309, 7).Value = "Ajo"
$ws.Cells.Item(309, 8).Value = "Chino"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 220
$ws.Cells.Item(309, 11).Value = 19000
$ws.Cells.Item(309, 12).Value = 20000
$ws.Cells.Item(309, 13).Value = 19455
$ws.Cells.Item(309, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(309, 15).Value = "China"
$ws.Cells.Item(309, 16).Value = 1946
$ws.Cells.Item(309, 17).Value = 10
$ws.Cells.Item(309, 18).Value = "Hortaliza"
